$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B19 was stored as text ("20650442421"); convert it to a real number to
# match the "All records" numeric-phone convention used elsewhere in the sheet.
$ws.Range("B19").Value = 20650442421

# Row 20: the validator's duplicate "Hina Kanjiani" failure, moved into the
# All-records sheet. Phone stored as a real number.
$ws.Range("A20").Value = "Hina Kanjiani"
$ws.Range("B20").Value = 20650442421
$ws.Range("C20").Value = "2050 Cross Creek Ct Allen 75013 TX USA"
$ws.Range("D20").Value = "GG"
$ws.Range("E20").Value = "English"
$ws.Range("K20").Value = "Phone validation failed: Invalid phone number format: 20650442421"
$ws.Range("L20").Value = "2025-09-15 16:36:55"
# Leading apostrophe forces text storage so "2025-09-15" isn't parsed into a date serial.
$ws.Range("M20").Formula = "'2025-09-15"
$ws.Range("N20").Value = "Failed"
$ws.Range("O20").Value = "Invalid Phone Number"

# Row 21: the validator's duplicate "Mohan Passi" failure, moved into the
# All-records sheet. Phone kept as text (leading apostrophe keeps the
# digit string from being re-typed as a number).
$ws.Range("A21").Value = "Mohan Passi"
$ws.Range("B21").Formula = "'20650442421"
$ws.Range("C21").Value = "4631 Gresham Drive, Eldorado Hills, 95762, CA, USA"
$ws.Range("D21").Value = "JKR"
$ws.Range("E21").Value = "English"
$ws.Range("K21").Value = "Phone validation failed: Invalid phone number format: 20650442421"
$ws.Range("L21").Value = "2025-09-15 16:37:06"
$ws.Range("M21").Formula = "'2025-09-15"
$ws.Range("N21").Value = "Failed"
$ws.Range("O21").Value = "Invalid Phone Number"
